$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.005723000769734084
$ws.Range("J2").Value = 0.005723000769734084
$ws.Range("M2").Value = 0.243056
$ws.Range("N2").Value = 0.729168
$ws.Range("O2").Value = 0.002199620488481675
$ws.Range("P2").Value = 0.002199620488481675
$ws.Range("Q2").Value = 0.06857938466133334
$ws.Range("R2").Value = 0.617214461952
$ws.Range("S2").Value = 0.00001258842974870349
$ws.Range("T2").Value = 0.00001258842974870349
$ws.Range("I3").Value = 0.005723000769734084
$ws.Range("J3").Value = 0.005723000769734084
$ws.Range("M3").Value = 70.95253000000001
$ws.Range("N3").Value = 212.85759
$ws.Range("O3").Value = 0.6421097964979703
$ws.Range("P3").Value = 0.6421097964979703
$ws.Range("Q3").Value = 20.01958745130667
$ws.Range("R3").Value = 180.17628706176
$ws.Range("S3").Value = 0.00367479485961168
$ws.Range("T3").Value = 0.00367479485961168
$ws.Range("I4").Value = 0.005723000769734084
$ws.Range("J4").Value = 0.005723000769734084
$ws.Range("M4").Value = 0.04794200000000001
$ws.Range("N4").Value = 0.143826
$ws.Range("O4").Value = 0.0004338679376719292
$ws.Range("P4").Value = 0.0004338679376719292
$ws.Range("Q4").Value = 0.01352705902933334
$ws.Range("R4").Value = 0.121743531264
$ws.Range("S4").Value = 0.000002483026541259391
$ws.Range("T4").Value = 0.000002483026541259391
$ws.Range("I5").Value = 0.005723000769734084
$ws.Range("J5").Value = 0.005723000769734084
$ws.Range("M5").Value = 39.25553366666666
$ws.Range("N5").Value = 117.766601
$ws.Range("O5").Value = 0.3552567150758761
$ws.Range("P5").Value = 0.3552567150758761
$ws.Range("Q5").Value = 11.07613201654044
$ws.Range("R5").Value = 99.685188148864
$ws.Range("S5").Value = 0.002033134453832441
$ws.Range("T5").Value = 0.002033134453832441
$ws.Range("G6").Value = 30.199365
$ws.Range("H6").Value = 90.598095
$ws.Range("I6").Value = 0.6125398923302606
$ws.Range("J6").Value = 0.6125398923302606
$ws.Range("M6").Value = 0.243056
$ws.Range("N6").Value = 0.729168
$ws.Range("O6").Value = 0.002199620488481675
$ws.Range("P6").Value = 0.002199620488481675
$ws.Range("Q6").Value = 7.340136859440001
$ws.Range("R6").Value = 66.06123173496
$ws.Range("S6").Value = 0.001347355297182
$ws.Range("T6").Value = 0.001347355297182
$ws.Range("G7").Value = 30.199365
$ws.Range("H7").Value = 90.598095
$ws.Range("I7").Value = 0.6125398923302606
$ws.Range("J7").Value = 0.6125398923302606
$ws.Range("M7").Value = 70.95253000000001
$ws.Range("N7").Value = 212.85759
$ws.Range("O7").Value = 0.6421097964979703
$ws.Range("P7").Value = 0.6421097964979703
$ws.Range("Q7").Value = 2142.72135114345
$ws.Range("R7").Value = 19284.49216029105
$ws.Range("S7").Value = 0.3933178656110723
$ws.Range("T7").Value = 0.3933178656110723
$ws.Range("G8").Value = 30.199365
$ws.Range("H8").Value = 90.598095
$ws.Range("I8").Value = 0.6125398923302606
$ws.Range("J8").Value = 0.6125398923302606
$ws.Range("M8").Value = 0.04794200000000001
$ws.Range("N8").Value = 0.143826
$ws.Range("O8").Value = 0.0004338679376719292
$ws.Range("P8").Value = 0.0004338679376719292
$ws.Range("Q8").Value = 1.44781795683
$ws.Range("R8").Value = 13.03036161147
$ws.Range("S8").Value = 0.0002657614198271158
$ws.Range("T8").Value = 0.0002657614198271158
$ws.Range("G9").Value = 30.199365
$ws.Range("H9").Value = 90.598095
$ws.Range("I9").Value = 0.6125398923302606
$ws.Range("J9").Value = 0.6125398923302606
$ws.Range("M9").Value = 39.25553366666666
$ws.Range("N9").Value = 117.766601
$ws.Range("O9").Value = 0.3552567150758761
$ws.Range("P9").Value = 0.3552567150758761
$ws.Range("Q9").Value = 1185.492189469455
$ws.Range("R9").Value = 10669.42970522509
$ws.Range("S9").Value = 0.2176089100021792
$ws.Range("T9").Value = 0.2176089100021792
$ws.Range("G10").Value = 18.820355
$ws.Range("H10").Value = 56.461065
$ws.Range("I10").Value = 0.3817371069000054
$ws.Range("J10").Value = 0.3817371069000054
$ws.Range("M10").Value = 0.243056
$ws.Range("N10").Value = 0.729168
$ws.Range("O10").Value = 0.002199620488481675
$ws.Range("P10").Value = 0.002199620488481675
$ws.Range("Q10").Value = 4.57440020488
$ws.Range("R10").Value = 41.16960184392
$ws.Range("S10").Value = 0.0008396767615509711
$ws.Range("T10").Value = 0.0008396767615509711
$ws.Range("G11").Value = 18.820355
$ws.Range("H11").Value = 56.461065
$ws.Range("I11").Value = 0.3817371069000054
$ws.Range("J11").Value = 0.3817371069000054
$ws.Range("M11").Value = 70.95253000000001
$ws.Range("N11").Value = 212.85759
$ws.Range("O11").Value = 0.6421097964979703
$ws.Range("P11").Value = 0.6421097964979703
$ws.Range("Q11").Value = 1335.35180274815
$ws.Range("R11").Value = 12018.16622473335
$ws.Range("S11").Value = 0.2451171360272864
$ws.Range("T11").Value = 0.2451171360272864
$ws.Range("G12").Value = 18.820355
$ws.Range("H12").Value = 56.461065
$ws.Range("I12").Value = 0.3817371069000054
$ws.Range("J12").Value = 0.3817371069000054
$ws.Range("M12").Value = 0.04794200000000001
$ws.Range("N12").Value = 0.143826
$ws.Range("O12").Value = 0.0004338679376719292
$ws.Range("P12").Value = 0.0004338679376719292
$ws.Range("Q12").Value = 0.90228545941
$ws.Range("R12").Value = 8.120569134690001
$ws.Range("S12").Value = 0.0001656234913035541
$ws.Range("T12").Value = 0.0001656234913035541
$ws.Range("G13").Value = 18.820355
$ws.Range("H13").Value = 56.461065
$ws.Range("I13").Value = 0.3817371069000054
$ws.Range("J13").Value = 0.3817371069000054
$ws.Range("M13").Value = 39.25553366666666
$ws.Range("N13").Value = 117.766601
$ws.Range("O13").Value = 0.3552567150758761
$ws.Range("P13").Value = 0.3552567150758761
$ws.Range("Q13").Value = 738.8030793211183
$ws.Range("R13").Value = 6649.227713890064
$ws.Range("S13").Value = 0.1356146706198645
$ws.Range("T13").Value = 0.1356146706198645
